$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" "27.013.37"
Set-TextValue "E2" "  -3.27%  "
Set-TextValue "D3" "1.727.53"
Set-TextValue "E3" "  -1.94%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "310.51"
Set-TextValue "E5" "  -5.36%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.18%  "
Set-TextValue "D7" "0.4827"
Set-TextValue "E7" "  +3.05%  "
Set-TextValue "D8" "0.3475"
Set-TextValue "E8" "  -1.01%  "
Set-TextValue "D9" "43.30"
Set-TextValue "E9" "  -0.40%  "
Set-TextValue "D10" "0.07229"
Set-TextValue "E10" "  -1.84%  "
Set-TextValue "D11" "1.051"
Set-TextValue "E11" "  -2.81%  "
Set-TextValue "E12" "  +0.23%  "
Set-TextValue "D13" "19.95"
Set-TextValue "E13" "  -2.91%  "
Set-TextValue "D14" "5.876"
Set-TextValue "E14" "  -1.93%  "
Set-TextValue "D15" "1.723.67"
Set-TextValue "E15" "  -2.02%  "
Set-TextValue "D16" "6.821"
Set-TextValue "E16" "  -4.61%  "
Set-TextValue "D17" "87.03"
Set-TextValue "E17" "  -5.84%  "
Set-TextValue "D18" "0.00001033"
Set-TextValue "E18" "  -1.87%  "
Set-TextValue "D19" "0.06400"
Set-TextValue "E19" "  -0.24%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  +0.09%  "
Set-TextValue "D21" "16.59"
Set-TextValue "E21" "  -1.30%  "
Set-TextValue "D22" "5.717"
Set-TextValue "E22" "  -0.65%  "
Set-TextValue "D23" "27.090.10"
Set-TextValue "E23" "  -3.02%  "
Set-TextValue "D24" "10.97"
Set-TextValue "E24" "  -1.47%  "
Set-TextValue "D25" "2.069"
Set-TextValue "E25" "  -3.76%  "
Set-TextValue "D26" "154.25"
Set-TextValue "E26" "  -4.90%  "
Set-TextValue "D27" "19.95"
Set-TextValue "E27" "  -0.21%  "
Set-TextValue "D28" "1.923.81"
Set-TextValue "E28" "  -1.86%  "
Set-TextValue "D29" "2.068"
Set-TextValue "E29" "  -4.50%  "
Set-TextValue "D30" "120.57"
Set-TextValue "E30" "  -1.86%  "
Set-TextValue "D31" "1.038"
Set-TextValue "E31" "  -3.01%  "
Set-TextValue "D32" "0.09306"
Set-TextValue "E32" "  -0.30%  "
Set-TextValue "D33" "3.642"
Set-TextValue "E33" "  -0.21%  "
Set-TextValue "D34" "5.379"
Set-TextValue "E34" "  -3.01%  "
Set-TextValue "D35" "0.05916"
Set-TextValue "E35" "  -2.41%  "
Set-TextValue "D36" "0.02184"
Set-TextValue "E36" "  -3.67%  "
Set-TextValue "D37" "1.431"
Set-TextValue "E37" "  +5.63%  "
Set-TextValue "B38" "Algorand"
Set-TextValue "C38" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2001"
Set-TextValue "E38" "  -3.06%  "
Set-TextValue "B39" "Aptos"
Set-TextValue "C39" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D39" "10.96"
Set-TextValue "E39" "  -5.84%  "
Set-TextValue "D40" "4.751"
Set-TextValue "E40" "  -2.99%  "
Set-TextValue "D41" "1.000"
Set-TextValue "E41" "  +0.23%  "
Set-TextValue "D42" "0.5985"
Set-TextValue "E42" "  -2.39%  "
Set-TextValue "D43" "1.116"
Set-TextValue "E43" "  -5.48%  "
Set-TextValue "D44" "7.535"
Set-TextValue "E44" "  -3.02%  "
Set-TextValue "D45" "12.74"
Set-TextValue "E45" "  -2.90%  "
Set-TextValue "E46" "  -4.13%  "
Set-TextValue "E47" "  -2.90%  "
Set-TextValue "D48" "119.32"
Set-TextValue "E48" "  -3.06%  "
Set-TextValue "E49" "  -3.81%  "
Set-TextValue "D50" "1.101"
Set-TextValue "E50" "  -1.81%  "
Set-TextValue "E51" "  -2.23%  "
